# UiComponentClassDiagram.pptx update
#  - updated DG (change app name, Component ClassDiagrams, Feature contribution,
#    Implementation of scheduling modules, tidy up user stories & NFR, add glossary)
#
# Concretely this applies to the pptx:
#  1. Refresh every "datetimeFigureOut" footer field (master / all layouts / notes
#     master) from 1/7/2017 -> 4/10/2018 (PowerPoint recaches these on save).
#  2. Rename the "BrowserPanel" class box on the UI component diagram to
#     "CalendarPanel".
#  3. Add a new connector/arrow shape (a duplicate of the existing "Freeform 115"
#     leader line) pointing at the (new) CalendarPanel box.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "1/7/2017") {
                $sh.TextFrame.TextRange.Text = "4/10/2018"
            }
        }
    }
}

# 1a. Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# 1b. Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# 1c. Notes master date placeholder.
Update-DatePlaceholder $p.NotesMaster.Shapes

# 2. Slide 1: rename BrowserPanel -> CalendarPanel.
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "BrowserPanel") {
            $sh.TextFrame.TextRange.Text = "CalendarPanel"
        }
    }
}

# 3. Add the new leader-line/arrow shape under the (renamed) CalendarPanel box.
#    It's a copy of the existing "Freeform 115" connector, flipped vertically,
#    repositioned/resized, with a solid (not dotted) dash.
$template = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.Name -eq "Freeform 115") {
        $template = $sh
    }
}

$dupRange = $template.Duplicate()
$newShape = $dupRange.Item(1)

$newShape.Left = 290.2488288976378
$newShape.Top = 303.28096488188976
$newShape.Width = 240.0
$newShape.Height = 25.066939133858266
$newShape.VerticalFlip = -1
